$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# New version-history rows to append after the last existing row (2.11 / 02/11/2018)
$rows = @(
    @("2.12", "Casos de prueba" + [char]0x2013 + "Se agregaron los casos de prueba.", "ntrezza", "04/11/2018"),
    @("2.13", "PERT " + [char]0x2013 + " Se realiz" + [char]0x00F3 + " un diagrama de pert.", "gbenitez", "04/11/2018"),
    @("2.14", "Casos de usos " + [char]0x2013 + " Se agregaron los escenarios de casos de uso y la especificacion de los casos de uso.", "gbenitez", "04/11/2018"),
    @("2.15", "Actualizacion en BD y crear aula " + [char]0x2013 + " Se realizaron algunos cambios en la base de datos y en la vista de crear aula.", "Jlucero - molmos", "05/11/2018"),
    @("2.16", "Microtaller " + [char]0x2013 + " Se agregaron las funcionalidades faltantes a la vista de microtalleres y se agrega una validacion en crear aula.", "molmos", "05/11/2018"),
    @("2.17", "Select dependientes " + [char]0x2013 + " Aplicaci" + [char]0x00F3 + "n de select dependientes a la hora de seleccionar departamentos. ", "molmos", "06/11/2018"),
    @("2.18", "Mas select dependientes y ajustes visuales " + [char]0x2013 + " Se colocaron los select dependientes que faltaban y se realizaron cambios esteticos solicitados por el cliente. Se mejoro la barra de navegaci" + [char]0x00F3 + "n.", "molmos", "07/11/2018"),
    @("2.19", "Video explicativo " + [char]0x2013 + " Se agrego un video en la vista de matricular usuarios.", "molmos", "08/11/2018"),
    @("2.20", "SideNav " + [char]0x2013 + " Se agreg" + [char]0x00F3 + " un Sidenav a la derecha en todas las vistas de la p" + [char]0x00E1 + "gina web.", "ntrezza", "08/11/2018"),
    @("2.21", "Arreglos en sidenav " + [char]0x2013 + " Se Arreglo el sidenav para que sea responsivo y se coloc" + [char]0x00F3 + " el logo.", "ntrezza - molmos", "08/11/2018")
)

foreach ($r in $rows) {
    $newRow = $tbl.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $r[0]
    $newRow.Cells.Item(2).Range.Text = $r[1]
    $newRow.Cells.Item(3).Range.Text = $r[2]
    $newRow.Cells.Item(4).Range.Text = $r[3]
}
